# Add two new columns, I ("I0") and J ("IF"), to the right of the
# existing H ("IP") column, with header formatting matching the other
# header cells and per-row numeric data for rows 2-45.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold font, border, centered) onto the
# new header cells I1/J1 before writing their text so the style index
# matches the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Per-row numeric values for the new I0 / IF columns, rows 2..45.
$i_vals = @(1,1,1,7,8,8,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,8,6,8,5,9,8,8,8,9,9,9,8,9,8,9,8,9,9,5,7,8,6)
$j_vals = @(1,5,5,7,8,8,5,5,6,5,6,5,4,5,5,6,4,4,5,6,5,4,8,6,8,7,9,9,8,8,9,9,9,8,9,9,9,9,9,9,7,9,8,6)

for ($r = 2; $r -le 45; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $i_vals[$idx]
    $ws.Cells.Item($r, 10).Value = $j_vals[$idx]
}
